# Applies cryptos.xlsx price/volume/row-order updates per commit
# "Updated cryptos list on Tue Mar 26 04:27:01 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.538.40'
$ws.Range('E2').Value = '  +5.14%  '
$ws.Range('D3').Value = '3.623.88'
$ws.Range('E3').Value = '  +5.05%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = "'591.74"
$ws.Range('E5').Value = '  +1.52%  '
$ws.Range('D6').Value = "'191.49"
$ws.Range('E6').Value = '  +3.97%  '
$ws.Range('D7').Value = "'0.646"
$ws.Range('E7').Value = '  +2.35%  '
$ws.Range('D8').Value = '3.615.82'
$ws.Range('E8').Value = '  +4.94%  '
$ws.Range('E10').Value = '  +2.75%  '
$ws.Range('D11').Value = "'0.666"
$ws.Range('E11').Value = '  +3.53%  '
$ws.Range('D12').Value = "'58.18"
$ws.Range('E12').Value = '  +3.42%  '
$ws.Range('D13').Value = "'0.0000289"
$ws.Range('E13').Value = '  +4.23%  '
$ws.Range('D14').Value = "'9.91"
$ws.Range('E14').Value = '  +5.37%  '
$ws.Range('D15').Value = '4.207.46'
$ws.Range('E15').Value = '  +5.03%  '
$ws.Range('D16').Value = "'19.73"
$ws.Range('E16').Value = '  +6.03%  '
$ws.Range('D17').Value = '3.625.54'
$ws.Range('D18').Value = '70.497.54'
$ws.Range('E18').Value = '  +4.94%  '
$ws.Range('D19').Value = "'12.66"
$ws.Range('E19').Value = '  +4.94%  '
$ws.Range('E20').Value = '  +0.50%  '
$ws.Range('E21').Value = '  +4.55%  '
$ws.Range('D22').Value = "'486.44"
$ws.Range('E22').Value = '  +1.22%  '
$ws.Range('D23').Value = "'19.67"
$ws.Range('E23').Value = '  +17.94%  '
$ws.Range('D24').Value = "'5.38"
$ws.Range('E24').Value = '  -1.80%  '
$ws.Range('D25').Value = "'4.45"
$ws.Range('E25').Value = '  +1.17%  '
$ws.Range('D26').Value = "'90.90"
$ws.Range('E26').Value = '  +1.66%  '
$ws.Range('D27').Value = "'3.14"
$ws.Range('E27').Value = '  +7.21%  '
$ws.Range('D28').Value = "'11.33"
$ws.Range('E28').Value = '  +3.34%  '
$ws.Range('D29').Value = "'9.66"
$ws.Range('E29').Value = '  +5.89%  '
$ws.Range('D30').Value = "'33.03"
$ws.Range('E30').Value = '  +5.56%  '
$ws.Range('D31').Value = "'7.84"
$ws.Range('E31').Value = '  +10.11%  '
$ws.Range('D32').Value = "'626.23"
$ws.Range('E32').Value = '  +6.56%  '
$ws.Range('D33').Value = "'12.29"
$ws.Range('E33').Value = '  +5.14%  '
$ws.Range('E34').Value = '  +7.86%  '
$ws.Range('D35').Value = "'66.44"
$ws.Range('E35').Value = '  +3.50%  '
$ws.Range('B36').Value = 'InjectiveProtocol'
$ws.Range('C36').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D36').Value = "'38.88"
$ws.Range('E36').Value = '  +6.83%  '
$ws.Range('B37').Value = 'TheGraph'
$ws.Range('C37').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D37').Value = "'0.412"
$ws.Range('E37').Value = '  +7.44%  '
$ws.Range('D38').Value = '0.0₃0815'
$ws.Range('E38').Value = '  +6.20%  '
$ws.Range('E39').Value = '  -1.34%  '
$ws.Range('E40').Value = '  -0.05%  '
$ws.Range('D41').Value = "'3.62"
$ws.Range('E41').Value = '  +1.91%  '
$ws.Range('D42').Value = '3.304.69'
$ws.Range('E42').Value = '  +3.12%  '
$ws.Range('E43').Value = '  +7.83%  '
$ws.Range('E44').Value = '  +10.82%  '
$ws.Range('E45').Value = '  +5.20%  '
$ws.Range('D46').Value = "'3.34"
$ws.Range('E46').Value = '  +4.28%  '
$ws.Range('E47').Value = '  +2.86%  '
$ws.Range('B48').Value = 'dogwifhat'
$ws.Range('C48').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D48').Value = "'2.74"
$ws.Range('E48').Value = '  +0.04%  '
$ws.Range('B49').Value = 'THORChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D49').Value = "'9.12"
$ws.Range('E49').Value = '  +4.49%  '
$ws.Range('D50').Value = "'3.31"
$ws.Range('E50').Value = '  +3.75%  '
$ws.Range('E51').Value = '  -0.04%  '
